# New test case row appended to the Leads test-case sheet (Sheet1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A25").Value = "CREATE LEAD 01"
$ws.Range("B25").Value = "click"
$ws.Range("C25").Value = "click"
$ws.Range("D25").Value = "no value"
$ws.Range("E25").Value = "null"
$ws.Range("F25").Value = "no"

# Matches the author's final selection/active-cell state after entering row 25.
$ws.Range("G25").Select()
